$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45954
$ws.Range("B2").Value = 101.85
$ws.Range("C2").Value = 92.51000000000001
$ws.Range("D2").Value = 88.63
$ws.Range("E2").Value = 87.64
$ws.Range("F2").Value = 87.84
$ws.Range("G2").Value = 90.73999999999999
$ws.Range("H2").Value = 105.9
$ws.Range("I2").Value = 134.67
$ws.Range("J2").Value = 160.68
$ws.Range("K2").Value = 133.3
$ws.Range("L2").Value = 93.69
$ws.Range("M2").Value = 62.51
$ws.Range("N2").Value = 48.32
$ws.Range("O2").Value = 47.44
$ws.Range("P2").Value = 47.01
$ws.Range("Q2").Value = 47.01
$ws.Range("R2").Value = 47.72
$ws.Range("S2").Value = 65.31999999999999
$ws.Range("T2").Value = 102.48
$ws.Range("U2").Value = 126.58
$ws.Range("V2").Value = 161.92
$ws.Range("W2").Value = 149.93
$ws.Range("X2").Value = 114.69
$ws.Range("Y2").Value = 103.96
$ws.Range("Z2").Value = 95.93000000000001
$ws.Range("AB2").Value = 132.62
$ws.Range("AD2").Value = 155.93
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 146.99
$ws.Range("AG2").Value = "1h-17h"
